$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: set a cell value as plain text, preventing Excel from
# auto-converting numeric-looking strings into numbers (which would
# drop meaningful trailing zeros), while leaving the cell style
# untouched (no explicit style index remains on the cell).
function Set-TextValue($cell, $val) {
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.Style = "Normal"
}

# Update Price (D) and Volume(1h) (E) columns for affected rows
Set-TextValue $ws.Cells.Item(2, 4) "56.587.41"
$ws.Cells.Item(2, 5).Value = "  -4.60%  "
Set-TextValue $ws.Cells.Item(3, 4) "2.353.14"
$ws.Cells.Item(3, 5).Value = "  -7.15%  "
$ws.Cells.Item(4, 5).Value = "  +0.10%  "
Set-TextValue $ws.Cells.Item(5, 4) "512.95"
$ws.Cells.Item(5, 5).Value = "  -4.12%  "
Set-TextValue $ws.Cells.Item(6, 4) "127.42"
$ws.Cells.Item(6, 5).Value = "  -6.68%  "
$ws.Cells.Item(7, 5).Value = "  -0.11%  "
$ws.Cells.Item(8, 5).Value = "  -3.28%  "
Set-TextValue $ws.Cells.Item(9, 4) "2.367.25"
$ws.Cells.Item(9, 5).Value = "  -6.57%  "
Set-TextValue $ws.Cells.Item(10, 4) "0.0950"
$ws.Cells.Item(10, 5).Value = "  -5.08%  "
$ws.Cells.Item(11, 5).Value = "  -2.24%  "
Set-TextValue $ws.Cells.Item(12, 4) "4.80"
$ws.Cells.Item(12, 5).Value = "  -8.28%  "
Set-TextValue $ws.Cells.Item(13, 4) "0.315"
$ws.Cells.Item(13, 5).Value = "  -6.66%  "
Set-TextValue $ws.Cells.Item(14, 4) "2.769.35"
$ws.Cells.Item(14, 5).Value = "  -7.18%  "
Set-TextValue $ws.Cells.Item(15, 4) "56.480.55"
$ws.Cells.Item(15, 5).Value = "  -4.74%  "
Set-TextValue $ws.Cells.Item(16, 4) "21.34"
$ws.Cells.Item(16, 5).Value = "  -5.64%  "
$ws.Cells.Item(17, 5).Value = "  -5.86%  "
Set-TextValue $ws.Cells.Item(18, 4) "2.336.31"
$ws.Cells.Item(18, 5).Value = "  -7.72%  "
Set-TextValue $ws.Cells.Item(19, 4) "10.26"
$ws.Cells.Item(19, 5).Value = "  -5.10%  "
Set-TextValue $ws.Cells.Item(20, 4) "308.67"
$ws.Cells.Item(20, 5).Value = "  -4.84%  "
$ws.Cells.Item(21, 5).Value = "  -6.24%  "
Set-TextValue $ws.Cells.Item(22, 4) "6.07"
$ws.Cells.Item(22, 5).Value = "  -1.30%  "
Set-TextValue $ws.Cells.Item(23, 4) "0.999"
$ws.Cells.Item(23, 5).Value = "  -0.09%  "
Set-TextValue $ws.Cells.Item(24, 4) "64.53"
$ws.Cells.Item(24, 5).Value = "  -2.07%  "
Set-TextValue $ws.Cells.Item(25, 4) "1.00"
$ws.Cells.Item(25, 5).Value = "  -0.02%  "
Set-TextValue $ws.Cells.Item(26, 4) "0.392"
$ws.Cells.Item(26, 5).Value = "  -4.62%  "
Set-TextValue $ws.Cells.Item(27, 4) "2.461.48"
$ws.Cells.Item(27, 5).Value = "  -7.06%  "
$ws.Cells.Item(28, 5).Value = "  -5.65%  "
Set-TextValue $ws.Cells.Item(29, 4) "7.14"
$ws.Cells.Item(29, 5).Value = "  -6.10%  "
Set-TextValue $ws.Cells.Item(30, 4) "171.94"
$ws.Cells.Item(30, 5).Value = "  -1.24%  "
$ws.Cells.Item(31, 5).Value = "  -5.53%  "
Set-TextValue $ws.Cells.Item(32, 4) "0.0₃0713"
$ws.Cells.Item(32, 5).Value = "  -7.48%  "
Set-TextValue $ws.Cells.Item(33, 4) "6.08"
$ws.Cells.Item(33, 5).Value = "  -5.43%  "
Set-TextValue $ws.Cells.Item(34, 4) "1.12"
$ws.Cells.Item(34, 5).Value = "  -8.44%  "
$ws.Cells.Item(35, 5).Value = "  -0.08%  "
$ws.Cells.Item(36, 5).Value = "  -0.25%  "
$ws.Cells.Item(37, 5).Value = "  -4.48%  "
Set-TextValue $ws.Cells.Item(38, 4) "1.18"
$ws.Cells.Item(38, 5).Value = "  -7.10%  "
Set-TextValue $ws.Cells.Item(39, 4) "3.69"
$ws.Cells.Item(39, 5).Value = "  -8.75%  "
Set-TextValue $ws.Cells.Item(42, 4) "1.42"
$ws.Cells.Item(42, 5).Value = "  -7.63%  "
Set-TextValue $ws.Cells.Item(43, 4) "3.30"
$ws.Cells.Item(43, 5).Value = "  -6.50%  "
Set-TextValue $ws.Cells.Item(44, 4) "4.85"
$ws.Cells.Item(44, 5).Value = "  -6.01%  "
Set-TextValue $ws.Cells.Item(45, 4) "123.20"
$ws.Cells.Item(45, 5).Value = "  -7.57%  "
Set-TextValue $ws.Cells.Item(46, 4) "0.567"
$ws.Cells.Item(46, 5).Value = "  -5.76%  "
Set-TextValue $ws.Cells.Item(47, 4) "251.12"
$ws.Cells.Item(47, 5).Value = "  -11.55%  "
$ws.Cells.Item(48, 5).Value = "  -2.96%  "
Set-TextValue $ws.Cells.Item(49, 4) "0.0486"
$ws.Cells.Item(49, 5).Value = "  -5.31%  "
$ws.Cells.Item(50, 5).Value = "  -6.77%  "
Set-TextValue $ws.Cells.Item(51, 4) "16.59"
$ws.Cells.Item(51, 5).Value = "  -7.34%  "

# Rows 40 and 41 swapped (SuiNetwork <-> OKB) with updated values
$ws.Cells.Item(40, 2).Value = "OKB"
$ws.Cells.Item(40, 3).Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
Set-TextValue $ws.Cells.Item(40, 4) "35.50"
$ws.Cells.Item(40, 5).Value = "  -3.48%  "

$ws.Cells.Item(41, 2).Value = "SuiNetwork"
$ws.Cells.Item(41, 3).Value = "https://coinranking.com/coin/3xJluUMvp+suinetwork-sui"
Set-TextValue $ws.Cells.Item(41, 4) "0.792"
$ws.Cells.Item(41, 5).Value = "  -0.33%  "
